$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.357.51'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.866.64'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.59'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.92%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4703'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.20'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.47%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.89'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.882.91'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6939'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.101'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '268.42'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.296.67'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.86'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007635'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.110.96'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.224'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.172'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.403'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.32'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.24%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.946'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.63%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.17%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.372'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.056'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04746'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.02%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7021'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.712'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.792'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.307'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.12'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.950'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4178'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.45%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8397'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.88'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '972.02'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.99%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.148'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.48'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05681'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.53%  '
